# Auto-generated PowerShell Excel COM-interop script
# Applies updated cryptocurrency price/volume values to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '60.910.96'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -3.18%  '
$ws.Range('E2').Style = "Normal"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.910.10'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -4.17%  '
$ws.Range('E3').Style = "Normal"
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('E4').Style = "Normal"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '586.16'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -1.40%  '
$ws.Range('E5').Style = "Normal"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '146.50'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -3.14%  '
$ws.Range('E6').Style = "Normal"
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E7').Style = "Normal"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.503'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -3.39%  '
$ws.Range('E8').Style = "Normal"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.909.32'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -4.05%  '
$ws.Range('E9').Style = "Normal"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.74'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +6.07%  '
$ws.Range('E10').Style = "Normal"
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -4.86%  '
$ws.Range('E11').Style = "Normal"
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -2.99%  '
$ws.Range('E12').Style = "Normal"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000224'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -4.22%  '
$ws.Range('E13').Style = "Normal"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '33.51'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -4.07%  '
$ws.Range('E14').Style = "Normal"
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.71%  '
$ws.Range('E15').Style = "Normal"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.394.71'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -4.16%  '
$ws.Range('E16').Style = "Normal"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '60.868.46'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -3.26%  '
$ws.Range('E17').Style = "Normal"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '6.76'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -3.85%  '
$ws.Range('E18').Style = "Normal"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '2.910.23'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -4.23%  '
$ws.Range('E19').Style = "Normal"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '426.59'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -6.42%  '
$ws.Range('E20').Style = "Normal"
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -4.53%  '
$ws.Range('E21').Style = "Normal"
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -3.02%  '
$ws.Range('E22').Style = "Normal"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.11'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -5.13%  '
$ws.Range('E23').Style = "Normal"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '80.17'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -3.40%  '
$ws.Range('E24').Style = "Normal"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '10.97'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('E25').Style = "Normal"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.20'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -2.98%  '
$ws.Range('E26').Style = "Normal"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.83'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -2.94%  '
$ws.Range('E27').Style = "Normal"
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('E28').Style = "Normal"
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('E29').Style = "Normal"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.25'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -2.18%  '
$ws.Range('E30').Style = "Normal"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.16'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +0.20%  '
$ws.Range('E32').Style = "Normal"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '26.47'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -4.29%  '
$ws.Range('E33').Style = "Normal"
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -4.02%  '
$ws.Range('E34').Style = "Normal"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0₃0844'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -1.24%  '
$ws.Range('E35').Style = "Normal"
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -2.52%  '
$ws.Range('E36').Style = "Normal"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '5.62'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -4.82%  '
$ws.Range('E37').Style = "Normal"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.97'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -5.01%  '
$ws.Range('E38').Style = "Normal"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '49.36'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -2.25%  '
$ws.Range('E39').Style = "Normal"
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -3.98%  '
$ws.Range('E40').Style = "Normal"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.124'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -0.64%  '
$ws.Range('E41').Style = "Normal"
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -4.91%  '
$ws.Range('E42').Style = "Normal"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.291'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -1.45%  '
$ws.Range('E43').Style = "Normal"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '41.58'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +0.93%  '
$ws.Range('E44').Style = "Normal"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '377.17'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -3.99%  '
$ws.Range('E45').Style = "Normal"
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -2.64%  '
$ws.Range('E46').Style = "Normal"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.671.43'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -2.72%  '
$ws.Range('E47').Style = "Normal"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '132.77'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('E48').Style = "Normal"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '24.73'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +2.92%  '
$ws.Range('E50').Style = "Normal"
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -1.75%  '
$ws.Range('E51').Style = "Normal"
